$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 247
$ws.Range("I5").Value = 90.5
$ws.Range("J5").Value = 481.75
$ws.Range("K5").Value = 90.5
$ws.Range("L5").Value = 481.75
$ws.Range("M5").Value = 24.5
$ws.Range("N5").Value = -711.75
$ws.Range("H19").Value = 5449.5
$ws.Range("I19").Value = 5449.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 5449.5
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -5274.5
$ws.Range("N19").ClearContents()
$ws.Range("H76").Value = 4998
$ws.Range("J76").Value = 4997
$ws.Range("L76").Value = 4997
$ws.Range("N76").Value = -5627
$ws.Range("H79").Value = 4998
$ws.Range("J79").Value = 4997
$ws.Range("L79").Value = 4997
$ws.Range("N79").Value = -7181
$ws.Range("H98").Value = 1952.3334
$ws.Range("I98").Value = 937
$ws.Range("J98").Value = 5506
$ws.Range("K98").Value = 937
$ws.Range("L98").Value = 5506
$ws.Range("M98").Value = 561
$ws.Range("N98").Value = -8502
$ws.Range("H122").Value = 1952.3334
$ws.Range("I122").Value = 937
$ws.Range("J122").Value = 5506
$ws.Range("K122").Value = 2811
$ws.Range("L122").Value = 16518
$ws.Range("M122").Value = -361
$ws.Range("N122").Value = -21418
$ws.Range("H138").Value = 3437.75
$ws.Range("J138").Value = 3133.6
$ws.Range("L138").Value = 9400.799999999999
$ws.Range("N138").Value = -19680.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1382.6923
$ws.Range("I2").Value = 618.1111
$ws.Range("K2").Value = 618.1111
$ws.Range("M2").Value = -505.1111
$ws.Range("H5").Value = 881
$ws.Range("J5").Value = 801.6667
$ws.Range("L5").Value = 801.6667
$ws.Range("N5").Value = -1025.6667
$ws.Range("H45").Value = 3877.7273
$ws.Range("I45").Value = 2486
$ws.Range("J45").Value = 5037.5
$ws.Range("K45").Value = 2486
$ws.Range("L45").Value = 5037.5
$ws.Range("M45").Value = -2109
$ws.Range("N45").Value = -5791.5
$ws.Range("H46").Value = 20983.334
$ws.Range("I46").Value = 19000
$ws.Range("J46").Value = 21975
$ws.Range("K46").Value = 19000
$ws.Range("L46").Value = 21975
$ws.Range("M46").Value = -18681
$ws.Range("N46").Value = -22613
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H116").Value = 1382.6923
$ws.Range("I116").Value = 618.1111
$ws.Range("K116").Value = 618.1111
$ws.Range("M116").Value = 1675.8889
$ws.Range("H134").Value = 64999
$ws.Range("J134").Value = 64999
$ws.Range("L134").Value = 64999
$ws.Range("N134").Value = -75139
$ws.Range("H135").Value = 6303749
$ws.Range("J135").Value = 61427.715
$ws.Range("L135").Value = 61427.715
$ws.Range("N135").Value = -71567.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1382.6923
$ws.Range("I3").Value = 618.1111
$ws.Range("K3").Value = 618.1111
$ws.Range("M3").Value = -504.1111
$ws.Range("H4").Value = 881
$ws.Range("J4").Value = 801.6667
$ws.Range("L4").Value = 801.6667
$ws.Range("N4").Value = -1031.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2638.9375
$ws.Range("I31").Value = 2051.7144
$ws.Range("K31").Value = 2051.7144
$ws.Range("M31").Value = -1756.7144
$ws.Range("H33").Value = 8739.799999999999
$ws.Range("J33").Value = 25999
$ws.Range("L33").Value = 25999
$ws.Range("N33").Value = -26757
$ws.Range("H34").Value = 2638.9375
$ws.Range("I34").Value = 2051.7144
$ws.Range("K34").Value = 2051.7144
$ws.Range("M34").Value = -1849.7144
$ws.Range("H99").Value = 4381.375
$ws.Range("I99").Value = 3699.7778
$ws.Range("K99").Value = 3699.7778
$ws.Range("M99").Value = -2201.7778
$ws.Range("H126").Value = 4381.375
$ws.Range("I126").Value = 3699.7778
$ws.Range("K126").Value = 11099.3334
$ws.Range("M126").Value = -8629.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 42197.145
$ws.Range("I22").Value = 47730
$ws.Range("K22").Value = 143190
$ws.Range("M22").Value = -143021
$ws.Range("H27").Value = 42197.145
$ws.Range("I27").Value = 47730
$ws.Range("K27").Value = 143190
$ws.Range("M27").Value = -143088
$ws.Range("H40").Value = 152.61539
$ws.Range("I40").Value = 125.6
$ws.Range("J40").Value = 242.66667
$ws.Range("K40").Value = 502.4
$ws.Range("L40").Value = 970.66668
$ws.Range("M40").Value = -433.4
$ws.Range("N40").Value = -1108.66668
$ws.Range("H92").Value = 876.41174
$ws.Range("I92").Value = 799.875
$ws.Range("K92").Value = 2399.625
$ws.Range("M92").Value = -1151.625
$ws.Range("H109").Value = 310
$ws.Range("I109").Value = 310
$ws.Range("K109").Value = 930
$ws.Range("M109").Value = 110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 2495
$ws.Range("I18").Value = 2000
$ws.Range("J18").Value = 2990
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 2990
$ws.Range("M18").Value = -1707
$ws.Range("N18").Value = -3576
$ws.Range("H70").Value = 11627.667
$ws.Range("I70").Value = 11627.667
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 11627.667
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -11357.667
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 11627.667
$ws.Range("I73").Value = 11627.667
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 11627.667
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -10691.667
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 1222
$ws.Range("I80").Value = 874
$ws.Range("J80").Value = 1361.2
$ws.Range("K80").Value = 874
$ws.Range("L80").Value = 1361.2
$ws.Range("M80").Value = 124
$ws.Range("N80").Value = -3357.2
$ws.Range("H83").Value = 1222
$ws.Range("I83").Value = 874
$ws.Range("J83").Value = 1361.2
$ws.Range("K83").Value = 4370
$ws.Range("L83").Value = 6806
$ws.Range("M83").Value = 622
$ws.Range("N83").Value = -16790
$ws.Range("H132").Value = 2539.625
$ws.Range("I132").Value = 2316.9546
$ws.Range("K132").Value = 6950.8638
$ws.Range("M132").Value = -4420.8638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2323.4707
$ws.Range("I7").Value = 1955
$ws.Range("J7").Value = 3521
$ws.Range("K7").Value = 1955
$ws.Range("L7").Value = 3521
$ws.Range("M7").Value = -1843
$ws.Range("N7").Value = -3745
$ws.Range("H14").Value = 7663
$ws.Range("I14").Value = 9999.5
$ws.Range("J14").Value = 2990
$ws.Range("K14").Value = 9999.5
$ws.Range("L14").Value = 2990
$ws.Range("M14").Value = -9827.5
$ws.Range("N14").Value = -3334
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H20").Value = 8333.333000000001
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -4774
$ws.Range("N20").Value = -15452
$ws.Range("H22").Value = 1858.0769
$ws.Range("I22").Value = 1144.375
$ws.Range("K22").Value = 1144.375
$ws.Range("M22").Value = -849.375
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 1858.0769
$ws.Range("I27").Value = 1144.375
$ws.Range("K27").Value = 1144.375
$ws.Range("M27").Value = -1037.375
$ws.Range("H46").Value = 3058.8235
$ws.Range("I46").Value = 2250
$ws.Range("K46").Value = 2250
$ws.Range("M46").Value = -2062
$ws.Range("H82").Value = 3374.125
$ws.Range("I82").Value = 798.3333
$ws.Range("K82").Value = 798.3333
$ws.Range("M82").Value = -437.3333
$ws.Range("H85").Value = 3374.125
$ws.Range("I85").Value = 798.3333
$ws.Range("K85").Value = 798.3333
$ws.Range("M85").Value = 449.6667
$ws.Range("H126").Value = 2323.4707
$ws.Range("I126").Value = 1955
$ws.Range("J126").Value = 3521
$ws.Range("K126").Value = 5865
$ws.Range("L126").Value = 10563
$ws.Range("M126").Value = -3395
$ws.Range("N126").Value = -15503
$ws.Range("H132").Value = 8938
$ws.Range("I132").Value = 8938
$ws.Range("K132").Value = 26814
$ws.Range("M132").Value = -24284
$ws.Range("H134").Value = 90000
$ws.Range("J134").Value = 90000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -100140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 55000
$ws.Range("J46").Value = 55000
$ws.Range("L46").Value = 55000
$ws.Range("N46").Value = -55462
$ws.Range("H134").Value = 55000
$ws.Range("J134").Value = 55000
$ws.Range("L134").Value = 165000
$ws.Range("N134").Value = -170070
